# Weekly fruit/vegetable price update: a new daily price record was
# inserted for "Haba" (Vega Modelo de Temuco) above the existing row 88,
# pushing the former rows 88-102 down to 89-103.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 88 - all rows from 88 downward
# shift down by one (88 -> 89, ..., 102 -> 103), carrying their values and
# formatting with them.
$ws.Rows.Item(88).EntireRow.Insert()

# Populate the newly inserted row 88 with the new price observation.
$ws.Cells.Item(88, 1).Value = 10
$ws.Cells.Item(88, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(88, 3).Value = "La Araucanía"
$ws.Cells.Item(88, 4).Value = 45244
$ws.Cells.Item(88, 5).Value = 9
$ws.Cells.Item(88, 6).Value = 100112026
$ws.Cells.Item(88, 7).Value = "Haba"
$ws.Cells.Item(88, 8).Value = "Sin especificar"
$ws.Cells.Item(88, 9).Value = "Primera"
$ws.Cells.Item(88, 10).Value = 50
$ws.Cells.Item(88, 11).Value = 10000
$ws.Cells.Item(88, 12).Value = 10000
$ws.Cells.Item(88, 13).Value = 10000
$ws.Cells.Item(88, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(88, 15).Value = "Región del Maule"
$ws.Cells.Item(88, 16).Value = 400
$ws.Cells.Item(88, 17).Value = 25
$ws.Cells.Item(88, 18).Value = "Hortaliza"
